# Modified framework with new URL
# Replace the contact-us sample data (subjectheading/email/orderref/message)
# with the new name/email/enquiry layout and drop column D entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: name / Gaurav / Dixit  (typed column-by-column, like the original author)
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = "Gaurav"
$ws.Range("A3").Value = "Dixit"

# Column B (email) keeps its existing values + hyperlinks, nothing to change there.

# Column C: enquiry / doing automation work / checking the automation
$ws.Range("C1").Value = "enquiry"
$ws.Range("C2").Value = "doing automation work"
$ws.Range("C3").Value = "checking the automation"

# Column D (orderref/message) is no longer used - clear it out entirely.
$ws.Range("D1:D3").Clear()

# Match the new active selection left behind in the saved worksheet.
$ws.Range("E8").Select() | Out-Null
